$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" contain the same refreshed data rows.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 9797
    $ws.Range("F4").Value = 37
    $ws.Range("F5").Value = 567
}
